# "Fruta / hortaliza, semanal" update
# A new weekly price observation is inserted at the top of the data
# (row 19, right after the sheet's frozen/earlier rows 2-18), pushing every
# existing record from row 19 downward by one row (old row 19 -> new row 20,
# ..., old row 90 -> new row 91). The sheet's used range grows from
# A1:R90 to A1:R91 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19; Excel shifts rows 19:90 down to 20:91
# and extends the used range accordingly.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new observation.
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44764
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 100112031
$ws.Range("G19").Value = "Poroto verde"
$ws.Range("H19").Value = "Magnum"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 32000
$ws.Range("L19").Value = 32000
$ws.Range("M19").Value = 32000
$ws.Range("N19").Value = "`$/malla 25 kilos"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 1280
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
